$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right above the current "Programa resumido:" row
# (row 12), pushing everything below it down by two rows.
$ws.Rows("12:13").Insert()

# New row 12 only carries the label in column A (bold style, like the
# other "header" rows e.g. "Avaliação:" / "Requisitos:").
$ws.Range("A12").Value = "Docentes responsáveis:"

# New row 13 carries the value, duplicated in columns B and C (wrapped /
# red-wrapped styles), matching the existing "label + duplicated value"
# pattern used throughout the sheet.
$ws.Range("B13:C13").Value = "8855158 - Morun Bernardino Neto"

# The row-insert operation copies formatting (and thus creates blank
# placeholder cells) across the whole A:C block from the row above.
# Remove the cells that should not exist on these two new rows so the
# row shape matches the rest of the sheet (label-only row / value-only
# row, as used elsewhere for single-column entries).
$ws.Range("B12:C12").Clear()
$ws.Range("A13").Clear()
